# Append the "Visual Drag-Drop Canvas Foundation Update" section
# (divider, heading, date, blank line, and two status-table rows) to the
# end of the document, mirroring the existing "Module Name / Developed /
# Partial Developed / Need To Develop" report blocks already in the file.

$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)  ; # wdCollapseEnd - move to the very end of the document

$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r></w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">---</w:t></w:r></w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Visual Drag-Drop Canvas Foundation Update</w:t></w:r></w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Updated: 2026-02-18</w:t></w:r></w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r></w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Module Name</w:t><w:tab/><w:t xml:space="preserve">Developed</w:t><w:tab/><w:t xml:space="preserve">Partial Developed</w:t><w:tab/><w:t xml:space="preserve">Need To Develop</w:t></w:r></w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Visual Drag-Drop Canvas</w:t><w:tab/><w:t xml:space="preserve">Added layout schema validation API, version history table/API, publish + rollback APIs, Store Builder undo/redo timeline, nested block (children) editing in section tree, responsive preview now renders nested blocks</w:t><w:tab/><w:t xml:space="preserve">WYSIWYG freeform drag surface is still simplified (list/tree interactions), nesting depth tooling basic, no collaborative editing</w:t><w:tab/><w:t xml:space="preserve">Full canvas interaction engine with true drag-n-drop between containers, deep nesting UX, advanced inspector panels, multi-user editing controls</w:t></w:r></w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Versioning &amp; Governance</w:t><w:tab/><w:t xml:space="preserve">Draft versions auto-created on save, publish/rollback wired to version IDs, version list shown in UI</w:t><w:tab/><w:t xml:space="preserve">No release notes per version, no lock/freeze workflow, no scheduled publish</w:t><w:tab/><w:t xml:space="preserve">Full version governance (approval gates, scheduled publish, diff viewer, immutable release snapshots)</w:t></w:r></w:p>
'@

$r.InsertXML($xml)

Write-Output "Appended Visual Drag-Drop Canvas Foundation Update section"
